# Update the "Gesamtinvestitionskosten" sheet with new input figures.
# Only the raw input cells (column B "netto" amounts, the manual D10
# tax-rate value and the B20 "% der Ust" rate) need to be written -
# all dependent formula cells (C, D, E, F columns, the sums in rows
# 12/14, and the cross-sheet references on the other worksheets)
# recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gesamtinvestitionskosten")

$ws.Range("B2").Value  = 100
$ws.Range("B3").Value  = 110
$ws.Range("B4").Value  = 120
$ws.Range("B5").Value  = 130
$ws.Range("B6").Value  = 140
$ws.Range("B7").Value  = 150
$ws.Range("B8").Value  = 160
$ws.Range("B9").Value  = 76
$ws.Range("B10").Value = 67

# Manually entered "% der Ust" value for the last row (not a formula).
$ws.Range("D10").Value = 0.77

# Shared Ust-rate used by the B3:B9 formulas ($B$20).
$ws.Range("B20").Value = 0.54
